# preparation publication 0.2.0
# - bump Version 0.1.1 -> 0.2.0
# - bump Date to 2023-10-20T08:59:58+00:00
# - add a new "Jurisdiction" / "iso:code:3166:FR" row on the Metadata sheet
#   (inserted after "Contact", pushing Description/Purpose/Copyright/Immutable
#   down by one row)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# --- Version & Date -------------------------------------------------------
$ws1.Range("B3").Value = "0.2.0"
$ws1.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# --- Insert "Jurisdiction" row at row 11 ----------------------------------
# Shift existing rows 11..14 (Description, Purpose, Copyright, Immutable)
# down to 12..15, working from the bottom up so nothing gets clobbered.
for ($r = 14; $r -ge 11; $r--) {
    $dest = $r + 1
    $ws1.Range("A$dest").Value = $ws1.Range("A$r").Value()
    $ws1.Range("B$dest").Value = $ws1.Range("B$r").Value()
}

# Row 15 is brand new territory on the sheet; give it the same formatting
# as its neighbour (row 14) before it was pushed down.
$ws1.Range("A14:B14").Copy()
$ws1.Range("A15:B15").PasteSpecial(-4122)

# Now populate the freed-up row 11 with the new property.
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = "iso:code:3166:FR"
